$wb = $excel.ActiveWorkbook

# --- Update the conversion text on sheet "Hoja1" (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$oldText = $cellA1.Value2
$newText = $oldText -replace [regex]::Escape("1000 Bs = 7.66 = 31053.91 pesos"), "1000 Bs = 7.41 = 30081.48 pesos"
$newText = $newText -replace [regex]::Escape("31053.91 pesos = 7.67 = 958.46 Bs"), "30081.48 pesos = 7.4 = 954.65 Bs"
$cellA1.Value = $newText

# --- Update the rate values on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 135
$wsTasas.Range("O10").Value = 4061
$wsTasas.Range("N12").Value = 4064.99
$wsTasas.Range("O12").Value = 129.005
